$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range('D2').Value = '57.043.11'
$ws.Range('E2').Value = '  -1.28%  '
$ws.Range('D3').Value = '2.981.82'
$ws.Range('E3').Value = '  -2.06%  '
Set-TextValue 'D4' '0.999'
$ws.Range('E4').Value = '  -0.06%  '
Set-TextValue 'D5' '501.92'
$ws.Range('E5').Value = '  -4.15%  '
Set-TextValue 'D6' '138.29'
$ws.Range('E6').Value = '  -3.06%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  -2.87%  '
Set-TextValue 'D9' '7.31'
$ws.Range('E9').Value = '  -4.09%  '
Set-TextValue 'D10' '0.108'
$ws.Range('E10').Value = '  -2.46%  '
$ws.Range('E11').Value = '  -1.94%  '
$ws.Range('D12').Value = '3.485.49'
$ws.Range('E12').Value = '  -2.18%  '
$ws.Range('E13').Value = '  -2.23%  '
Set-TextValue 'D14' '26.04'
$ws.Range('E14').Value = '  -0.92%  '
Set-TextValue 'D15' '0.0000160'
$ws.Range('E15').Value = '  -2.00%  '
$ws.Range('D16').Value = '57.068.67'
$ws.Range('E16').Value = '  -1.19%  '
$ws.Range('E17').Value = '  -1.61%  '
$ws.Range('D18').Value = '2.985.28'
$ws.Range('E18').Value = '  -2.15%  '
$ws.Range('E19').Value = '  -2.28%  '
$ws.Range('E20').Value = '  -2.76%  '
Set-TextValue 'D21' '320.93'
$ws.Range('E21').Value = '  -5.05%  '
$ws.Range('E22').Value = '  -0.21%  '
Set-TextValue 'D23' '5.72'
$ws.Range('E23').Value = '  +0.24%  '
$ws.Range('E24').Value = '  -1.53%  '
Set-TextValue 'D25' '63.80'
$ws.Range('E25').Value = '  -1.81%  '
Set-TextValue 'D26' '0.999'
$ws.Range('E26').Value = '  +0.24%  '
Set-TextValue 'D27' '0.165'
$ws.Range('E27').Value = '  -5.35%  '
$ws.Range('D28').Value = '0.0₃0898'
$ws.Range('E28').Value = '  -5.34%  '
Set-TextValue 'D29' '6.56'
$ws.Range('E29').Value = '  -5.01%  '
Set-TextValue 'D30' '7.07'
$ws.Range('E30').Value = '  -2.37%  '
$ws.Range('E31').Value = '  -4.03%  '
$ws.Range('E32').Value = '  -5.01%  '
Set-TextValue 'D33' '20.17'
$ws.Range('E33').Value = '  -3.65%  '
Set-TextValue 'D34' '155.48'
$ws.Range('E34').Value = '  -1.80%  '
Set-TextValue 'D35' '4.58'
$ws.Range('E35').Value = '  -2.81%  '
$ws.Range('E36').Value = '  -1.52%  '
$ws.Range('E37').Value = '  -5.10%  '
Set-TextValue 'D38' '24.06'
$ws.Range('E38').Value = '  -4.04%  '
Set-TextValue 'D39' '0.0667'
$ws.Range('E39').Value = '  -3.51%  '
Set-TextValue 'D40' '37.83'
$ws.Range('E40').Value = '  +0.50%  '
$ws.Range('D41').Value = '3.011.37'
$ws.Range('E41').Value = '  -2.06%  '
Set-TextValue 'D42' '0.999'
$ws.Range('E42').Value = '  -0.17%  '
Set-TextValue 'D43' '3.74'
$ws.Range('E43').Value = '  -1.82%  '
Set-TextValue 'D44' '0.640'
$ws.Range('E44').Value = '  -2.97%  '
$ws.Range('E45').Value = '  -4.63%  '
$ws.Range('D46').Value = '2.199.10'
$ws.Range('E46').Value = '  -5.97%  '
Set-TextValue 'D47' '0.947'
$ws.Range('E47').Value = '  -7.60%  '
$ws.Range('E48').Value = '  -0.53%  '
Set-TextValue 'D49' '0.0236'
$ws.Range('E49').Value = '  -4.41%  '
Set-TextValue 'D50' '19.21'
$ws.Range('E50').Value = '  -3.33%  '
$ws.Range('E51').Value = '  -10.36%  '
